$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$text = "Conversión del día 💰`n" +
        "✅ Dólar paralelo: 68`n" +
        "`n" +
        "Binance`n" +
        "✅ 1000 Bs = 1.81 = 6537.22 pesos`n" +
        "✅ 6537.22 pesos = 1.8 = 924.91 Bs`n" +
        "`n" +
        "Promedio competencia`n" +
        "✅ Tasa pesos: 20`n" +
        "✅ Tasa Bs: 20`n" +
        "✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $text

# --- tasas: update the N10/O10/N12/O12 rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 553.14
$ws2.Range("O10").Value = 3616
$ws2.Range("N12").Value = 3640
$ws2.Range("O12").Value = 515.001
